# Brooklyn roster: rotate the three rows covering Dru Smith (TW), Mikal
# Bridges and Cameron Johnson so the players line up with the correct
# jersey #, height, weight, birth date, experience and college — the
# shared-string table in the source file had drifted out of sync with the
# row data for rows 15-17.
#
# Net effect (row "No." stays put, only the player record shifts):
#   row 15 (No. 13) : Dru Smith (TW)  -> Mikal Bridges
#   row 16 (No. 14) : Mikal Bridges   -> Cameron Johnson
#   row 17 (No. 15) : Cameron Johnson -> Dru Smith (TW)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "Mikal Bridges"
$ws.Range("D15").Value = "SF"
$ws.Range("E15").Value = "6-6"
$ws.Range("F15").Value = 209
$ws.Range("G15").Value = "August 30, 1996"
$ws.Range("H15").Value = "us"
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "4"
$ws.Range("I15").ClearFormats()
$ws.Range("J15").Value = "Villanova"
$ws.Range("K15").Value = "https://www.basketball-reference.com/players/b/bridgmi01.html"

$ws.Range("B16").Value = 2
$ws.Range("C16").Value = "Cameron Johnson"
$ws.Range("D16").Value = "PF"
$ws.Range("E16").Value = "6-8"
$ws.Range("F16").Value = 210
$ws.Range("G16").Value = "March 3, 1996"
$ws.Range("H16").Value = "us"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "3"
$ws.Range("I16").ClearFormats()
$ws.Range("J16").Value = "Pitt, UNC"
$ws.Range("K16").Value = "https://www.basketball-reference.com/players/j/johnsca02.html"

$ws.Range("B17").Value = 9
$ws.Range("C17").Value = "Dru Smith (TW)"
$ws.Range("D17").Value = "SG"
$ws.Range("E17").Value = "6-3"
$ws.Range("F17").Value = 203
$ws.Range("G17").Value = "December 30, 1997"
$ws.Range("H17").Value = "us"
$ws.Range("I17").Value = "R"
$ws.Range("J17").Value = "University of Evansville, Missouri"
$ws.Range("K17").Value = "https://www.basketball-reference.com/players/s/smithdr01.html"
